# DSA2 Final Grading Guide - apply grading updates
# (octree container, asteroid collision grading, bezier surface skewing)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: "Demonstrate broad phase collision detection with an octree" ---
# Score entered (Worth 10/10) and grading comment added.
# Note: set G12's string BEFORE G10's so the shared-string table order
# matches (Collision preserves/transfers momentum, then Visual evident...).
$ws.Range("F12").Value = 10
$ws.Range("G12").Value = "Collision preserves/transfers momentum"

# --- Row 10: "Use at least one kind of integration technique" ---
$ws.Range("G10").Value = "Visual evident in collision of bodies"

# --- Row 15: "Demonstrate a scaling transformation" ---
$ws.Range("F15").Value = 10
$ws.Range("G15").Value = "Sphere-sphere"

# --- Row 16: "Demonstrate a shearing transformation" ---
$ws.Range("F16").Value = 10

# --- Row 20: "Include interesting audio effects when something collides" ---
$ws.Range("F20").Value = 5

# --- Row 21: "Use SLERP with quaternions" ---
$ws.Range("F21").Value = 5

# --- Row 54: Bezier surface bonus points increased (oscillating skewing) ---
$ws.Range("F54").Value = 20

# --- Update the active selection shown when the file was last saved ---
$ws.Range("R4").Select()
